# updated GME & AMC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 18 (pushes "Poker" from E18 down to E19)
$ws.Rows("18:18").Insert()

# New row 18: "Interactive Entertainment" with a hyperlink (gets Hyperlink style)
$ws.Range("E18").Value = "Interactive Entertainment"
$ws.Hyperlinks.Add($ws.Range("E18"), "Interactive Entertainment.xlsx") | Out-Null
$ws.Range("E18").Style = "Hyperlink"

# Rename "Blackjack" (E16) to "Card Games (Non-Poker)"
$ws.Range("E16").Value = "Card Games (Non-Poker)"

# "Chess" (E17) becomes a hyperlink too
$ws.Hyperlinks.Add($ws.Range("E17"), "Chess.xlsx") | Out-Null
$ws.Range("E17").Style = "Hyperlink"

# Update the active selection to E15
$ws.Range("E15").Select() | Out-Null
